# Update the CL max (takeoff) input values on Sheet1.
# These are the raw input cells; all other cells in the workbook
# (including the derived formulas on Sheet1 and the other sheets,
# as well as the chart caches) recalculate automatically from them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 2.7
$ws.Range("G2").Value = 1.9
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 3.2
$ws.Range("G4").Value = 2.5

$excel.CalculateFullRebuild()

foreach ($ws2 in $wb.Worksheets) {
    foreach ($co in $ws2.ChartObjects()) {
        $co.Chart.Refresh()
    }
}
